# Applies the "Changelog 25 th April 2023" edits:
#  1. Merge the "A survey on a sample..." runs (the proofErr-bounded
#     ") ," run disappears) into one contiguous run of text.
#  2. Merge the "Find the number of cars..." runs (the proofErr-bounded
#     "W  (" run disappears) into one contiguous run of text.
#  3. Insert the missing "A " after "Grade " in the 50-students question,
#     landing as its own run (split out of the original run) the way
#     Word leaves things after an in-place correction.

$d = $word.ActiveDocument

# --- Edit 1: consolidate the "survey on a sample of 25 new cars" runs --
$old1 = "A survey on a sample of 25 new cars being sold at a local auto dealer was conducted to see which of three popular options, air conditioning (A), radio (R), and power windows (W) , were already installed. The survey found"
$new1 = "A survey on a sample of 25 new cars being sold at a local auto dealer was conducted to see which of three popular options, air conditioning (A), radio (R), and power windows (W) , were already installed. The survey found"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Edit 2: consolidate the "Find the number of cars..." runs --------
$old2 = "Find the number of cars that had (a) only W  (b) only A, (c) only R, (d) R and W but not A (e) A and R but nor W  (f) only one of the option (g) at least one option, (k) None of the option. "
$new2 = "Find the number of cars that had (a) only W  (b) only A, (c) only R, (d) R and W but not A (e) A and R but nor W  (f) only one of the option (g) at least one option, (k) None of the option. "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Edit 3: insert "A " after "Grade " (grade-A wording fix) ----------
# First drop in a unique marker right where the missing word belongs,
# then overwrite the marker (bolding/un-bolding it) so the insertion
# survives as a standalone run instead of re-merging with its neighbours.
$d.Content.Find.Execute("who got Grade in exactly", $true, $false, $false, $false, $false, $true, 1, $false, "who got Grade |TMPMARK|in exactly", 2) | Out-Null

$markRange = $d.Content
$markRange.Find.Execute("|TMPMARK|", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markRange.Bold = 1
$markRange.Text = "A "

$wordRange = $d.Content
$wordRange.Find.Execute("Grade A in exactly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordRange.MoveStart(1, 6) | Out-Null
$wordRange.MoveEnd(1, -10) | Out-Null
$wordRange.Bold = 0
